$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: new contact "De716"
$ws.Range("A18").Value = "De716"
$ws.Range("B18").Value = "tinoveler@gmail.com"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = "0650591900"

# Row 19: new contact "COT"
$ws.Range("A19").Value = "COT"
$ws.Range("B19").Value = "enaccf31400@gmail.com"
$ws.Hyperlinks.Add($ws.Range("B19"), "mailto:enaccf31400@gmail.com") | Out-Null
$ws.Range("B19").Style = "Lien hypertexte"
$ws.Range("C19").Value = "0650591900"

# Update the selected range shown in the sheet view
$ws.Range("F24").Select() | Out-Null
